$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column C
$ws.Range("C609").Value = 2
$ws.Range("C966").Value = 10
$ws.Range("C967").Value = 32
$ws.Range("C971").Value = 31
$ws.Range("C975").Value = 34
$ws.Range("C979").Value = 16
$ws.Range("C980").Value = 34
$ws.Range("C987").Value = 32
$ws.Range("C990").Value = 6
$ws.Range("C991").Value = 13
$ws.Range("C992").Value = 27
$ws.Range("C993").Value = 3
$ws.Range("C997").Value = 35
$ws.Range("C999").Value = 8
$ws.Range("C1000").Value = 27
$ws.Range("C1001").Value = 26
$ws.Range("C1005").Value = 10
$ws.Range("C1006").Value = 10
$ws.Range("C1007").Value = 21
$ws.Range("C1008").Value = 2
$ws.Range("C1009").Value = 5
$ws.Range("C1010").Value = 17
$ws.Range("C1011").Value = 35
$ws.Range("C1014").Value = 5
$ws.Range("C1015").Value = 15
$ws.Range("C1016").Value = 31

# Add new rows for date 44166 (2020-12-01)
$ws.Range("A1017").Value = 44166
$ws.Range("B1017").Value = "50-59"
$ws.Range("C1017").Value = 1

$ws.Range("A1018").Value = 44166
$ws.Range("B1018").Value = "60-69"
$ws.Range("C1018").Value = 2

$ws.Range("A1019").Value = 44166
$ws.Range("B1019").Value = "70-79"
$ws.Range("C1019").Value = 8

$ws.Range("A1020").Value = 44166
$ws.Range("B1020").Value = "80+"
$ws.Range("C1020").Value = 17

# Copy style (including date number format) from the previous date row to new rows
$ws.Range("A1016").Copy() | Out-Null
$ws.Range("A1017:A1020").PasteSpecial(-4122) | Out-Null
